$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date serials from 45190 to 45192 for rows 2-11
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
